$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new measurement (2026/02/19, 木, 0, 36) was recorded between the existing
# 2026/02/18 and 2026/12/29 entries. Insert a new row at 813, shifting the
# old rows 813:854 down to 814:855, then fill in the new row's data.
$ws.Rows.Item(813).Insert([Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown)

# Force column A to be treated as plain text so the date-like string isn't
# auto-converted into a date serial number/format (matching the sibling
# rows, which store the date as a literal string).
$ws.Cells.Item(813, 1).NumberFormat = "@"
$ws.Cells.Item(813, 1).Value = "2026/02/19"
$ws.Cells.Item(813, 1).ClearFormats()

$ws.Cells.Item(813, 2).Value = "木"
$ws.Cells.Item(813, 3).Value = 0
$ws.Cells.Item(813, 4).Value = 36
